$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates (weekly data shift) ---
$ws.Range("D2").Value = 44330
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("S2").Value = 861

# --- Row 3 is unchanged ---

# --- Row 4 updates ---
$ws.Range("D4").Value = 44698
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16500
$ws.Range("Q4").Value = "$/caja 18 kilos granel"
$ws.Range("S4").Value = 917
$ws.Range("T4").Value = 18

# --- Row 5 updates ---
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 17500
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17750
$ws.Range("S5").Value = 1109

# --- Row 6 updates ---
$ws.Range("D6").Value = 44316
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("Q6").Value = "$/caja 16 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 16

# --- New row 7 (appended record) ---
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 44344
$ws.Range("D7").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = "Otros"
$ws.Range("I7").Value = 100107001
$ws.Range("J7").Value = "Caqui"
$ws.Range("K7").Value = "Mankaki"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 13500
$ws.Range("Q7").Value = "$/caja 18 kilos granel"
$ws.Range("R7").Value = "Provincia de Curicó"
$ws.Range("S7").Value = 750
$ws.Range("T7").Value = 18
